$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New target values scraped from the NATMI re-run ("Natmi following Dr Hou advice").
# Keys are A1-style cell references on the single data sheet; values are the
# literal numbers the diff shows for the `after` OOXML.
$updates = @{
    "E2" = 3
    "G2" = 23.061728
    "H2" = 69.18518399999999
    "I2" = 0.6130144106248721
    "J2" = 0.6902769593117909
    "K2" = 3
    "M2" = 40.91730166666667
    "N2" = 122.751905
    "O2" = 0.2897771170516138
    "P2" = 0.3083463959441224
    "Q2" = 943.6236815306133
    "R2" = 8492.613133775518
    "S2" = 0.1776375486219696
    "T2" = 0.2128444126070583
    "E3" = 3
    "G3" = 23.061728
    "H3" = 69.18518399999999
    "I3" = 0.6130144106248721
    "J3" = 0.6902769593117909
    "K3" = 3
    "M3" = 34.738136
    "N3" = 104.214408
    "O3" = 0.2460161470038337
    "P3" = 0.2617811683839066
    "Q3" = 801.1214436590079
    "R3" = 7210.092992931071
    "S3" = 0.150811443359757
    "T3" = 0.180701508917131
    "E4" = 3
    "G4" = 23.061728
    "H4" = 69.18518399999999
    "I4" = 0.6130144106248721
    "J4" = 0.6902769593117909
    "K4" = 3
    "M4" = 18.806737
    "N4" = 56.42021099999999
    "O4" = 0.1331896729995656
    "P4" = 0.1417246332776418
    "Q4" = 433.7158532615359
    "R4" = 3903.442679353823
    "S4" = 0.08164718889514816
    "T4" = 0.09782924891846924
    "E5" = 3
    "G5" = 23.061728
    "H5" = 69.18518399999999
    "I5" = 0.6130144106248721
    "J5" = 0.6902769593117909
    "K5" = 3
    "M5" = 21.229913
    "N5" = 63.689739
    "O5" = 0.1503506520179033
    "P5" = 0.1599853092240957
    "Q5" = 489.598479069664
    "R5" = 4406.386311626976
    "S5" = 0.09216711633382026
    "T5" = 0.1104341727857654
    "E6" = 3
    "G6" = 23.061728
    "H6" = 69.18518399999999
    "I6" = 0.6130144106248721
    "J6" = 0.6902769593117909
    "K6" = 2
    "M6" = 25.510579
    "N6" = 51.021158
    "O6" = 0.1806664109270835
    "P6" = 0.1281624931702333
    "Q6" = 588.318034020512
    "R6" = 3529.908204123072
    "S6" = 0.1107511134141771
    "T6" = 0.08846761608336683
    "E7" = 3
    "G7" = 1.926013333333334
    "H7" = 5.778040000000001
    "I7" = 0.05119624723650278
    "J7" = 0.05764887294340218
    "K7" = 3
    "M7" = 40.91730166666667
    "N7" = 122.751905
    "O7" = 0.2897771170516138
    "P7" = 0.3083463959441224
    "Q7" = 78.80726857402223
    "R7" = 709.2654171662001
    "S7" = 0.01483550092805542
    "T7" = 0.01777582220233869
    "E8" = 3
    "G8" = 1.926013333333334
    "H8" = 5.778040000000001
    "I8" = 0.05119624723650278
    "J8" = 0.05764887294340218
    "K8" = 3
    "M8" = 34.738136
    "N8" = 104.214408
    "O8" = 0.2460161470038337
    "P8" = 0.2617811683839066
    "Q8" = 66.90611311114667
    "R8" = 602.1550180003201
    "S8" = 0.01259510348618008
    "T8" = 0.0150913893151392
    "E9" = 3
    "G9" = 1.926013333333334
    "H9" = 5.778040000000001
    "I9" = 0.05119624723650278
    "J9" = 0.05764887294340218
    "K9" = 3
    "M9" = 18.806737
    "N9" = 56.42021099999999
    "O9" = 0.1331896729995656
    "P9" = 0.1417246332776418
    "Q9" = 36.22202621849333
    "R9" = 325.99823596644
    "S9" = 0.00681881142823472
    "T9" = 0.008170265376773041
    "E10" = 3
    "G10" = 1.926013333333334
    "H10" = 5.778040000000001
    "I10" = 0.05119624723650278
    "J10" = 0.05764887294340218
    "K10" = 3
    "M10" = 21.229913
    "N10" = 63.689739
    "O10" = 0.1503506520179033
    "P10" = 0.1599853092240957
    "Q10" = 40.88909550350667
    "R10" = 368.0018595315601
    "S10" = 0.007697389152877975
    "T10" = 0.009222972764270799
    "E11" = 3
    "G11" = 1.926013333333334
    "H11" = 5.778040000000001
    "I11" = 0.05119624723650278
    "J11" = 0.05764887294340218
    "K11" = 2
    "M11" = 25.510579
    "N11" = 51.021158
    "O11" = 0.1806664109270835
    "P11" = 0.1281624931702333
    "Q11" = 49.13371529505334
    "R11" = 294.80229177032
    "S11" = 0.009249442241154576
    "T11" = 0.00738842328488043
    "E12" = 2
    "G12" = 12.632464
    "H12" = 25.264928
    "I12" = 0.3357893421386252
    "J12" = 0.2520741677448068
    "K12" = 3
    "M12" = 40.91730166666667
    "N12" = 122.751905
    "O12" = 0.2897771170516138
    "P12" = 0.3083463959441224
    "Q12" = 516.8863402813066
    "R12" = 3101.318041687839
    "S12" = 0.09730406750158879
    "T12" = 0.07772616113472533
    "E13" = 2
    "G13" = 12.632464
    "H13" = 25.264928
    "I13" = 0.3357893421386252
    "J13" = 0.2520741677448068
    "K13" = 3
    "M13" = 34.738136
    "N13" = 104.214408
    "O13" = 0.2460161470038337
    "P13" = 0.2617811683839066
    "Q13" = 438.8282524471039
    "R13" = 2632.969514682623
    "S13" = 0.08260960015789665
    "T13" = 0.06598827015163641
    "E14" = 2
    "G14" = 12.632464
    "H14" = 25.264928
    "I14" = 0.3357893421386252
    "J14" = 0.2520741677448068
    "K14" = 3
    "M14" = 18.806737
    "N14" = 56.42021099999999
    "O14" = 0.1331896729995656
    "P14" = 0.1417246332776418
    "Q14" = 237.575428109968
    "R14" = 1425.452568659808
    "S14" = 0.04472367267618275
    "T14" = 0.03572511898239952
    "E15" = 2
    "G15" = 12.632464
    "H15" = 25.264928
    "I15" = 0.3357893421386252
    "J15" = 0.2520741677448068
    "K15" = 3
    "M15" = 21.229913
    "N15" = 63.689739
    "O15" = 0.1503506520179033
    "P15" = 0.1599853092240957
    "Q15" = 268.186111695632
    "R15" = 1609.116670173792
    "S15" = 0.05048614653120513
    "T15" = 0.04032816367405948
    "E16" = 2
    "G16" = 12.632464
    "H16" = 25.264928
    "I16" = 0.3357893421386252
    "J16" = 0.2520741677448068
    "K16" = 2
    "M16" = 25.510579
    "N16" = 51.021158
    "O16" = 0.1806664109270835
    "P16" = 0.1281624931702333
    "Q16" = 322.261470836656
    "R16" = 1289.045883346624
    "S16" = 0.06066585527175192
    "T16" = 0.03230645380198605
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

Write-Output ("Updated " + $updates.Count + " cells")
